{"js": "// Resume edits:\n// 1. Summary paragraph: drop \"junior \" from \"A motivated junior developer...\"\n// 2. Summary paragraph: the lone \". \" run right after \"...learning how to be\n//    precise.\" loses its period (becomes a single space) since the next\n//    sentence run gets folded together with it.\n// 3. Work-experience paragraph (\"Served as a technical instructor...\"):\n//    - \"lesson online\" -> \"lessons in person and\"\n//    - \"Piloted new grading systems m\" -> \"Piloted and mastered new grading\n//      systems and programs of the institution. Brought an open and\n//      appealing environment to each class. \"\n\nconst body = context.document.body;\n\nconst edits = [\n  {\n    find: \"A motivated junior developer with 2+ years of experience\",\n    replace: \"A motivated developer with 2+ years of experience\",\n  },\n  {\n    find: \"learning how to be precise.. \",\n    replace: \"learning how to be precise. \",\n  },\n  {\n    find: \"lesson online\",\n    replace: \"lessons in person and\",\n  },\n  {\n    find: \"Piloted new grading systems m\",\n    replace:\n      \"Piloted and mastered new grading systems and programs of the institution. Brought an open and appealing environment to each class. \",\n  },\n];\n\nfor (const { find, replace } of edits) {\n  const results = body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + find);\n  }\n\n  results.items[0].insertText(replace, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Resume edits:\n# 1. Summary paragraph: drop \"junior \" from \"A motivated junior developer...\"\n# 2. Summary paragraph: the lone \". \" run right after \"...learning how to be\n#    precise.\" loses its period (becomes a single space) since the next\n#    sentence run gets folded together with it.\n# 3. Work-experience paragraph (\"Served as a technical instructor...\"):\n#    - \"lesson online\" -> \"lessons in person and\"\n#    - \"Piloted new grading systems m\" -> \"Piloted and mastered new grading\n#      systems and programs of the institution. Brought an open and\n#      appealing environment to each class. \"\n\n$d = $word.ActiveDocument\n\nfunction Replace-DocText($findText, $replaceWith) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceWith\n\n    # Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n    #   ReplaceWith, Replace)\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $found = $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceWith, 2)\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\nReplace-DocText \"A motivated junior developer with 2+ years of experience\" \"A motivated developer with 2+ years of experience\"\nReplace-DocText \"learning how to be precise.. \" \"learning how to be precise. \"\nReplace-DocText \"lesson online\" \"lessons in person and\"\nReplace-DocText \"Piloted new grading systems m\" \"Piloted and mastered new grading systems and programs of the institution. Brought an open and appealing environment to each class. \"\n"}
